$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the "Conversión del día" rate text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$text = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.56 = 9620.61 pesos`n✅ 9620.61 pesos = 2.55 = 949.43 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $text

# --- tasas sheet: refreshed exchange-rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 389.89
$ws2.Range("O10").Value = 3750.98
$ws2.Range("N12").Value = 3769.6
$ws2.Range("O12").Value = 372.009
